# The commit swaps the two theme parts in the package:
#   ppt/theme/theme1.xml  ("Integral" / "Red Violet")   -> becomes the stock "Office Theme" palette
#   ppt/theme/theme2.xml  ("Office Theme" / "Office")   -> becomes the "Integral" / "Red Violet" palette
#
# theme1.xml is the theme that is actually wired to the deck's single slide master
# (ppt/slideMasters/slideMaster1.xml -> ppt/theme/theme1.xml), so it is what is reachable
# from the PowerPoint object model via SlideMaster.Theme.ThemeColorScheme. The font scheme
# (majorFont/minorFont) and the format scheme (fill/line/effect styles) are byte-for-byte
# identical between the two themes already, so only the 12 color-scheme slots need to move.
#
# Re-point every slot of the slide master's theme color scheme to the stock Office palette.
# PowerPoint's RGB() helper isn't available in this host, so the long values below are
# written directly as 0xBBGGRR literals (COM long colour order) - the RRGGBB value each
# one represents is noted in the trailing comment.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1      -> RRGGBB 000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      -> RRGGBB FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      -> RRGGBB 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      -> RRGGBB E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  -> RRGGBB 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  -> RRGGBB ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  -> RRGGBB A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  -> RRGGBB FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  -> RRGGBB 4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  -> RRGGBB 70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    -> RRGGBB 0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink -> RRGGBB 954F72
